# CU04 Realizar Llamada.xlsx - fix references to "CU Cambio Estado Oportunidad"
# so that they correctly point to "CU06 Cambio Estado Oportunidad".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generar Ticket")

# Paso 6 (Puntos de Extensión) - row 12
$ws.Range("B12").Value = "Paso 6 - El sistema invoca CU06 Cambio Estado Oportunidad"

# Curso Básico, paso 6 - row 19
$ws.Range("B19").Value = "El sistema modifica el estado de la oportunidad. Invocando CU06 Cambio Estado Oportunidad"

# Update view: remove the frozen/scrolled topLeftCell and move the active selection to B19
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B19").Select()
